$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting the rest of the table (old rows 6-31)
# down to become rows 7-32. This makes room for the new
# "05_03_export_bootstraps.sas" line that documents the new export step.
$ws.Rows.Item(6).Insert()

# Row 5 ("05_02_bootstrap.sas") keeps its program/description text but its
# output-file cell is updated, and the row shrinks to fit the shorter text.
$ws.Range("C5").Value = "flows_jtw1990_moe.sas7bdat bootclusters_jtw1990_moe.sas7bdat"
$ws.Rows.Item(5).RowHeight = 35.05

# New row 6 documents the new export/stats step.
$ws.Range("A6").Value = "05_03_export_bootstraps.sas"
$ws.Range("B6").Value = "Exports data and does some stats"
$ws.Range("C6").Value = "bootclusters_jtw1990_moe.csv, bootclusters_jtw1990_moe_new.dta, flows_jtw1990_moe.dta, flows_jtw1990_moe.csv"
$ws.Rows.Item(6).RowHeight = 46.25

# The active selection moves from B8 to C8.
$ws.Range("C8").Select()
